$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Date column to Text format so date-like strings are not
# auto-converted to date serial numbers when assigned below.
$ws.Range("BF2:BF31").NumberFormat = "@"

$ws.Range("BF2").Value = "2014-02-02"
$ws.Range("AA3").Value = 19
$ws.Range("AC3").Value = -4.3
$ws.Range("AD3").Value = 5
$ws.Range("AE3").Value = 26
$ws.Range("AG3").Value = 27
$ws.Range("AH3").Value = 28
$ws.Range("AM3").Value = 24
$ws.Range("AO3").Value = 23
$ws.Range("AP3").Value = 25
$ws.Range("AQ3").Value = 12
$ws.Range("AU3").Value = 29
$ws.Range("AX3").Value = 16
$ws.Range("BA3").Value = 28
$ws.Range("BF3").Value = "2014-02-02"
$ws.Range("D3").Value = 48
$ws.Range("E3").Value = 15
$ws.Range("G3").Value = 0.313
$ws.Range("I3").Value = 36
$ws.Range("K3").Value = 0.434
$ws.Range("L3").Value = 6.4
$ws.Range("M3").Value = 19.5
$ws.Range("N3").Value = 0.326
$ws.Range("O3").Value = 16.2
$ws.Range("P3").Value = 21.2
$ws.Range("Q3").Value = 0.767
$ws.Range("V3").Value = 15.6
$ws.Range("X3").Value = 4.7
$ws.Range("Z3").Value = 21.7
$ws.Range("AQ4").Value = 11
$ws.Range("BF4").Value = "2014-02-02"
$ws.Range("BF5").Value = "2014-02-02"
$ws.Range("BF6").Value = "2014-02-02"
$ws.Range("AO7").Value = 24
$ws.Range("AU7").Value = 28
$ws.Range("BF7").Value = "2014-02-02"
$ws.Range("AD8").Value = 5
$ws.Range("BF8").Value = "2014-02-02"
$ws.Range("AV9").Value = 16
$ws.Range("BF9").Value = "2014-02-02"
$ws.Range("AM10").Value = 25
$ws.Range("BF10").Value = "2014-02-02"
$ws.Range("AD11").Value = 5
$ws.Range("AO11").Value = 25
$ws.Range("AP11").Value = 20
$ws.Range("BF11").Value = "2014-02-02"
$ws.Range("BF12").Value = "2014-02-02"
$ws.Range("BF13").Value = "2014-02-02"
$ws.Range("BF14").Value = "2014-02-02"
$ws.Range("BF15").Value = "2014-02-02"
$ws.Range("BA16").Value = 29
$ws.Range("BF16").Value = "2014-02-02"
$ws.Range("AZ17").Value = 8
$ws.Range("BF17").Value = "2014-02-02"
$ws.Range("BF18").Value = "2014-02-02"
$ws.Range("BF19").Value = "2014-02-02"
$ws.Range("BF20").Value = "2014-02-02"
$ws.Range("BF21").Value = "2014-02-02"
$ws.Range("BF22").Value = "2014-02-02"
$ws.Range("AB23").Value = 96.40000000000001
$ws.Range("AD23").Value = 5
$ws.Range("AH23").Value = 5
$ws.Range("AM23").Value = 17
$ws.Range("AV23").Value = 17
$ws.Range("AZ23").Value = 9
$ws.Range("BF23").Value = "2014-02-02"
$ws.Range("D23").Value = 48
$ws.Range("F23").Value = 35
$ws.Range("G23").Value = 0.271
$ws.Range("I23").Value = 36.3
$ws.Range("J23").Value = 82.2
$ws.Range("K23").Value = 0.442
$ws.Range("M23").Value = 20.6
$ws.Range("Q23").Value = 0.758
$ws.Range("S23").Value = 33.2
$ws.Range("T23").Value = 42.3
$ws.Range("Y23").Value = 6
$ws.Range("Z23").Value = 20.1
$ws.Range("AD24").Value = 5
$ws.Range("AE24").Value = 26
$ws.Range("AG24").Value = 27
$ws.Range("AH24").Value = 6
$ws.Range("BF24").Value = "2014-02-02"
$ws.Range("BF25").Value = "2014-02-02"
$ws.Range("AX26").Value = 18
$ws.Range("BF26").Value = "2014-02-02"
$ws.Range("AE27").Value = 26
$ws.Range("AG27").Value = 26
$ws.Range("AU27").Value = 27
$ws.Range("BF27").Value = "2014-02-02"
$ws.Range("BF28").Value = "2014-02-02"
$ws.Range("BF29").Value = "2014-02-02"
$ws.Range("AX30").Value = 17
$ws.Range("BF30").Value = "2014-02-02"
$ws.Range("AM31").Value = 18
$ws.Range("AP31").Value = 26
$ws.Range("BF31").Value = "2014-02-02"
